$d = $word.ActiveDocument

# --- 1. Add the new paragraph styles (based on MSCJoin) ---
$sA = $d.Styles.Add("MSC_Join_A", 1)
$sA.BaseStyle = $d.Styles("MSCJoin")

$sB = $d.Styles.Add("MSC_Join_B", 1)
$sB.BaseStyle = $d.Styles("MSCJoin")
$sB.Font.Name = "Noto Sans CJK SC"
$sB.Font.NameFarEast = "Noto Sans CJK SC"
$sB.Font.NameBi = "Noto Sans CJK SC"

$sC = $d.Styles.Add("MSC_Join_C", 1)
$sC.BaseStyle = $d.Styles("MSCJoin")

# --- 2. Re-point the "MSCJoin" paragraphs that sit between an
#        "MSC_Paragraph_A" block and the next "MSC_Paragraph_A"/"B" block
#        onto the freshly minted styles, matching the English/Chinese
#        verse groups in the document. ---

$paras = $d.Paragraphs
$count = $paras.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Style.NameLocal -eq "MSC_Join") {
        # Determine which "column" (A/B) this join paragraph belongs to by
        # looking back to the nearest preceding MSC_Paragraph_A / MSC_Paragraph_B.
        $j = $i - 1
        $target = $null
        while ($j -ge 1) {
            $prevStyle = $paras.Item($j).Style.NameLocal
            if ($prevStyle -eq "MSC_Paragraph_A") {
                $target = "MSC_Join_A"
                break
            }
            if ($prevStyle -eq "MSC_Paragraph_B") {
                $target = "MSC_Join_B"
                break
            }
            if ($prevStyle -eq "MSC_Paragraph_C") {
                break
            }
            $j = $j - 1
        }
        if ($target -ne $null) {
            $p.Style = $d.Styles($target)
        }
    }
}
